$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker list (rows 1-501) reflecting the "ABD dip hourly update"
$tickers = @(
    'Sembol',
    'MSFT',
    'AVGO',
    'TSLA',
    'BRK-A',
    'BRK-B',
    'JPM',
    'V',
    'ORCL',
    'PLTR',
    'NFLX',
    'JPM-PD',
    'JPM-PC',
    'UNH',
    'BAC-PK',
    'GE',
    'BAC-PL',
    'SAP',
    'BAC-PE',
    'BAC-PB',
    'AXP',
    'RY',
    'TMO',
    'CRM',
    'KLAC',
    'WFC-PY',
    'DIS',
    'BA',
    'WFC-PL',
    'APH',
    'ISRG',
    'ABT',
    'BX',
    'APP',
    'SHOP',
    'ACN',
    'UBER',
    'DHR',
    'TJX',
    'BKNG',
    'QCOM',
    'SPGI',
    'INTU',
    'PDD',
    'NOW',
    'COF',
    'BSX',
    'NEM',
    'SONY',
    'PANW',
    'ADBE',
    'PGR',
    'ARM',
    'CRWD',
    'KKR',
    'SBUX',
    'AEM',
    'CVS',
    'BN',
    'SPOT',
    'ADP',
    'CEG',
    'CVNA',
    'GD',
    'SNPS',
    'HOOD',
    'NKE',
    'MCO',
    'WM',
    'DASH',
    'CDNS',
    'B',
    'HWM',
    'NTES',
    'ELV',
    'MMM',
    'MS-PK',
    'AMT',
    'BAM',
    'MS-PI',
    'ABNB',
    'REGN',
    'MS-PF',
    'MS-PE',
    'CTAS',
    'DELL',
    'APO',
    'CI',
    'SE',
    'SNOW',
    'RELX',
    'MRVL',
    'WBD',
    'WPM',
    'TEL',
    'NET',
    'AJG',
    'NSC',
    'CNI',
    'DUK-PA',
    'RACE',
    'RKT',
    'CTA-PB',
    'URI',
    'ADSK',
    'COIN',
    'VLO',
    'NDAQ',
    'SRE',
    'IDXX',
    'TRI',
    'ZTS',
    'AU',
    'BIDU',
    'PYPL',
    'VST',
    'CMG',
    'RBLX',
    'EA',
    'SCHW-PD',
    'WDAY',
    'FNV',
    'EW',
    'GFI',
    'ARES',
    'CRWV',
    'PSA',
    'AXON',
    'DDOG',
    'ALNY',
    'CTA-PA',
    'HEI',
    'MSTR',
    'KGC',
    'TTWO',
    'SPG-PJ',
    'EXC',
    'XEL',
    'WCN',
    'ROP',
    'JD',
    'DHI',
    'RKLB',
    'KR',
    'MET-PA',
    'CTSH',
    'TCOM',
    'LVS',
    'AIG',
    'IQV',
    'RDDT',
    'CPRT',
    'XYZ',
    'HEI-A',
    'ALC',
    'A',
    'PAYX',
    'CCI',
    'KDP',
    'MDLN',
    'RMD',
    'FICO',
    'VEEV',
    'GEHC',
    'FISV',
    'TEAM',
    'CPNG',
    'OTIS',
    'CLS',
    'SYM',
    'APO-PA',
    'UAL',
    'ZS',
    'INSM',
    'NTRA',
    'KMB',
    'RJF',
    'MDB',
    'PCG',
    'EXPE',
    'CHT',
    'UI',
    'DG',
    'PSA-PK',
    'SOFI',
    'ALL-PB',
    'FOXA',
    'HUM',
    'WTW',
    'EXR',
    'FIS',
    'FOX',
    'VRSK',
    'QSR',
    'FLUT',
    'VICI',
    'AMRZ',
    'MTD',
    'TSCO',
    'NRG',
    'SYF',
    'HPE',
    'TME',
    'CBOE',
    'DTE',
    'STLA',
    'KHC',
    'CSGP',
    'PHG',
    'IRM',
    'ALAB',
    'PAAS',
    'WIT',
    'LEN',
    'ATO',
    'BRO',
    'EXE',
    'NMR',
    'LEN-B',
    'EFX',
    'ES',
    'TPR',
    'FSLR',
    'MKL',
    'DLTR',
    'AWK',
    'OMC',
    'WRB',
    'AVB',
    'DLR-PK',
    'RGLD',
    'BR',
    'LDOS',
    'PSTG',
    'SQM',
    'AXIA-PC',
    'ILMN',
    'TPG',
    'WAT',
    'TW',
    'CRDO',
    'NBIS',
    'CPAY',
    'LULU',
    'DRI',
    'OWL',
    'CNC',
    'IP',
    'FUTU',
    'AFRM',
    'DLR-PJ',
    'TLK',
    'CYBR',
    'FWONK',
    'ALB',
    'PSLV',
    'FWONA',
    'CG',
    'RL',
    'CMS',
    'UTHR',
    'GPN',
    'AS',
    'SSNC',
    'TWLO',
    'Q',
    'INCY',
    'HL',
    'SBAC',
    'WWD',
    'CHKP',
    'PTC',
    'NTAP',
    'TOST',
    'SGI',
    'GIB',
    'RIVN',
    'WY',
    'PODD',
    'TYL',
    'IFF',
    'RVMD',
    'KTOS',
    'HIG-PG',
    'GRAB',
    'IOT',
    'FTV',
    'PHYS',
    'U',
    'SMCI',
    'DKS',
    'HPQ',
    'CRCL',
    'LII',
    'FITBI',
    'XPEV',
    'IT',
    'PSNYW',
    'AGI',
    'ALLY',
    'PINS',
    'ESS',
    'SN',
    'CRS',
    'WST',
    'HUBS',
    'ZBH',
    'LNT',
    'LI',
    'ZG',
    'ATI',
    'TRMB',
    'MEDP',
    'QXO',
    'TXT',
    'Z',
    'HOLX',
    'THC',
    'APTV',
    'TKO',
    'CDW',
    'WES',
    'TRU',
    'TTD',
    'INVH',
    'MKC',
    'CDE',
    'NLY',
    'J',
    'MKC-V',
    'MAA',
    'OKTA',
    'COO',
    'TLN',
    'GFL',
    'WMG',
    'KSPI',
    'ROKU',
    'ROIV',
    'GEN',
    'H',
    'HTHT',
    'IONQ',
    'ONON',
    'DKNG',
    'AA',
    'AVAV',
    'MGA',
    'KEY-PI',
    'EMA',
    'CSL',
    'HMY',
    'ULS',
    'GDDY',
    'ARCC',
    'ICLR',
    'FIG',
    'RGC',
    'KRMN',
    'MAS',
    'RF-PC',
    'W',
    'CACI',
    'JHX',
    'EVR',
    'CELH',
    'DPZ',
    'BBY',
    'SBSW',
    'LOGI',
    'GWRE',
    'SOLV',
    'RBRK',
    'FIGR',
    'NLY-PG',
    'BLD',
    'NLY-PF',
    'AMH',
    'OHI',
    'RVTY',
    'JKHY',
    'RYAN',
    'BLDR',
    'SF',
    'PSKY',
    'UHS',
    'CHWY',
    'UNM',
    'SNAP',
    'JEF',
    'BNT',
    'BJ',
    'OKLO',
    'HLI',
    'GLPI',
    'EQX',
    'IVZ',
    'ACM',
    'BF-B',
    'SWK',
    'AGNC',
    'BMNR',
    'GLXY',
    'AG',
    'GMED',
    'AMKR',
    'IAG',
    'DT',
    'ZBRA',
    'TXRH',
    'DOC',
    'SMMT',
    'EPAM',
    'ALGN',
    'JOBY',
    'SUZ',
    'CR',
    'CMA',
    'TEM',
    'WYNN',
    'EXEL',
    'NTNX',
    'BSY',
    'NIO',
    'RDY',
    'DOCU',
    'BXP',
    'CEF',
    'MDGL',
    'BAH',
    'SCI',
    'MICC',
    'CNM',
    'MP',
    'WTRG',
    'DAY',
    'SARO',
    'CRL',
    'AR',
    'MOH',
    'MANH',
    'FDS',
    'AFG',
    'CART',
    'AIT',
    'PAG',
    'TECH',
    'CAE',
    'YMM',
    'ENSG',
    'PCOR',
    'KLAR',
    'ARMK',
    'JAZZ',
    'CHYM',
    'SAIL',
    'SANM',
    'BVN',
    'NGD',
    'POOL',
    'BIO-B',
    'EHC',
    'ALV',
    'SOLS',
    'ARE',
    'EGO',
    'TAP',
    'DDS',
    'BROS',
    'AAL',
    'ORI',
    'AYI',
    'ABVX',
    'ACI',
    'HSIC',
    'KNSL',
    'IESC',
    'MGM',
    'UWMC',
    'GS-PC',
    'RGEN',
    'DOX',
    'STEP',
    'LINE',
    'SWKS',
    'MORN',
    'UEC',
    'QBTS',
    'TTAN',
    'RZB',
    'AGNCM',
    'AGNCN',
    'LKQ',
    'PEGA',
    'CUBE',
    'VIPS',
    'OR',
    'AUR',
    'PL',
    'LEVI',
    'ADC'
)

$n = $tickers.Count
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $tickers[$i]
}

$targetRange = "A1:A" + $n
$ws.Range($targetRange).Value = $arr

# Remove the now-stale tail rows (previously A1:A521) so the used range shrinks to A1:A501
$ws.Range("A502:A521").ClearContents()

